# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.789.78"
$ws.Range("E2").Value = "  +5.52%  "

$ws.Range("D3").Value = "2.228.93"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.70"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.35%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.28"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0899"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.78%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "2.560.80"
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.03"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.803"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.59"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").Value = "2.244.00"
$ws.Range("E18").Value = "  +3.56%  "

$ws.Range("D19").Value = "41.718.23"
$ws.Range("E19").Value = "  +5.32%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.10"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.97%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.142"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.91"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.79%  "

$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.39%  "

$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0636"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.09%  "

$ws.Range("E37").Value = "  -4.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.66"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.91%  "

$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000254"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +28.93%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0239"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.41%  "

$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.58"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0976"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.12%  "

$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.55%  "

$ws.Range("D48").Value = "1.482.18"
$ws.Range("E48").Value = "  -2.10%  "

$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.48"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.33%  "
